# login() is added with negative testcases
# Reconstructs the new "COD Order" row (row 5), re-flows the old rows down,
# adds a "MailingReport" row and a new "Testing Error page" row, per the
# TestSuit.xlsx diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 loses its trailing columns (F:I) entirely ---------------------
$ws.Range("F6:I6").Clear()

# --- Cell values -----------------------------------------------------------
$ws.Range("C3").Value  = "YES"

$ws.Range("A5").Value  = "COD Order"
$ws.Range("B5").Value  = "“5”"
$ws.Range("C5").Value  = "YES"
$ws.Range("D5").Value  = "login"
$ws.Range("E5").Value  = "clearCart"
$ws.Range("F5").Value  = "productCatalogPage"
$ws.Range("G5").Value  = "productDetailPage"
$ws.Range("H5").Value  = "cartCheck"
$ws.Range("I5").Value  = "checkout"
$ws.Range("J5").Value  = "orderCOD"

$ws.Range("A6").Value  = "Checking Filters"
$ws.Range("B6").Value  = "“4”"
$ws.Range("C6").Value  = "YES"
$ws.Range("D6").Value  = "productCatalogPage"
$ws.Range("E6").Value  = "applyFilters"

$ws.Range("A7").Value  = "MailingReport"
$ws.Range("B7").Value  = "“”"
$ws.Range("D7").Value  = "emailReport"
$ws.Range("E7").ClearContents()

$ws.Range("A8").Value  = "Testing Error page"
$ws.Range("B8").Value  = "“6”"
$ws.Range("C8").Value  = "NO"
$ws.Range("D8").Value  = "tesetErrorpage"
$ws.Range("E8").Value  = "tesetErrorpage"

# --- Styles: copy formats from cells that already carry the right style ---
# style 0 = default, style 1 = bold header, style 2 = "NO/YES" style,
# style 3 = Action1-ish monospace, style 4 = Action2-ish monospace
$ws.Range("A3").Copy()
$ws.Range("C5,H5,I5,J5,A8,B8,C8,D8,E8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A2").Copy()
$ws.Range("G5,C6,E6,C7,E7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("D2").Copy()
$ws.Range("D5,E5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("E3").Copy()
$ws.Range("F5,D6,D7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = $false
